{"js": "// Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Z GC tests\n//\n// The document is a single-column table where every row holds one\n// benchmark statistic. This script:\n//   1. Updates the heap-size summary rows (rows 0-3).\n//   2. Removes three stale GC-pause rows and refreshes the remaining\n//      ones with the corrected values.\n//   3. Inserts three additional GC-pause rows with the real figures.\n//   4. Collapses the three tab-separated \"everything on one line\"\n//      summary rows near the end into single clean values.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Replace the text of a single-paragraph cell while keeping the\n// existing run formatting (rFonts/sz) intact, by replacing the\n// paragraph's range contents instead of the whole cell body.\nasync function setCellText(rowIndex, text) {\n  const cell = table.getCell(rowIndex, 0);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  const para = paragraphs.items[0];\n  const rng = para.getRange();\n  rng.insertText(text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- rows 0-2: \"100\"/\"0\"/\"47\" -> \"0M\" ---\nawait setCellText(0, \"0M\");\nawait setCellText(1, \"0M\");\nawait setCellText(2, \"0M\");\n\n// --- row 3: \"3\" -> \"42\" ---\nawait setCellText(3, \"42\");\n\n// --- delete the three now-obsolete rows that followed row 3\n//     (previously \"0.00003\", \"0.00006\", \"0.00004\") ---\ntable.rows.load(\"items\");\nawait context.sync();\ntable.rows.items[6].delete();\nawait context.sync();\ntable.rows.items[5].delete();\nawait context.sync();\ntable.rows.items[4].delete();\nawait context.sync();\n\n// After the deletions above, the rows shift up: the row that used to\n// be at index 7 is now at index 4 (unchanged \"0.00002\"), and indices\n// 5-8 hold the values that need refreshed text.\nawait setCellText(5, \"0.00009\");\nawait setCellText(6, \"0.00005\");\nawait setCellText(7, \"0.00001\");\nawait setCellText(8, \"0.00008\");\n\n// --- insert three brand-new rows right after row 8 ---\ntable.rows.load(\"items\");\nawait context.sync();\nconst anchorRow = table.rows.items[8];\nanchorRow.insertRows(Word.InsertLocation.after, 3, [\n  [\"0.00008\"],\n  [\"0.00009\"],\n  [\"0.00170\"],\n]);\nawait context.sync();\n\n// --- collapse the three multi-run/tab-separated summary rows near\n//     the end of the table into single plain values ---\nawait setCellText(43, \"100\");\nawait setCellText(44, \"0\");\nawait setCellText(45, \"47\");\n", "ps1": "# Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Z GC tests\n#\n# The document is a single-column table where every row holds one\n# benchmark statistic. This script:\n#   1. Updates the heap-size summary rows (1-based rows 1-4).\n#   2. Removes three stale GC-pause rows and refreshes the remaining\n#      ones with the corrected values.\n#   3. Inserts three additional GC-pause rows with the real figures.\n#   4. Collapses the three tab-separated \"everything on one line\"\n#      summary rows near the end into single clean values.\n#\n# Note: Word's Table/Rows/Cell collections are 1-based.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# --- rows 1-3 (1-based): \"100\"/\"0\"/\"47\" -> \"0M\" ---\n$tbl.Cell(1, 1).Range.Text = \"0M\"\n$tbl.Cell(2, 1).Range.Text = \"0M\"\n$tbl.Cell(3, 1).Range.Text = \"0M\"\n\n# --- row 4 (1-based): \"3\" -> \"42\" ---\n$tbl.Cell(4, 1).Range.Text = \"42\"\n\n# --- delete the three now-obsolete rows that followed row 4\n#     (previously \"0.00003\", \"0.00006\", \"0.00004\"; 1-based items 5,6,7).\n#     Delete from the bottom up so the remaining indices don't shift\n#     under us. ---\n$tbl.Rows.Item(7).Delete()\n$tbl.Rows.Item(6).Delete()\n$tbl.Rows.Item(5).Delete()\n\n# After the deletions above, the rows shift up: the row that used to\n# be 1-based item 8 is now item 5 (unchanged \"0.00002\"), and items\n# 6-9 hold the values that need refreshed text.\n$tbl.Cell(6, 1).Range.Text = \"0.00009\"\n$tbl.Cell(7, 1).Range.Text = \"0.00005\"\n$tbl.Cell(8, 1).Range.Text = \"0.00001\"\n$tbl.Cell(9, 1).Range.Text = \"0.00008\"\n\n# --- insert three brand-new rows right after 1-based row 9.\n#     Rows.Add(ref) always inserts immediately before the reference\n#     row, so use the following row (1-based item 10) as the anchor\n#     and add the new rows in reverse order to end up with the\n#     desired top-to-bottom sequence. ---\n$refRow = $tbl.Rows.Item(10)\n$newRow3 = $tbl.Rows.Add($refRow)\n$newRow3.Cells.Item(1).Range.Text = \"0.00170\"\n$newRow2 = $tbl.Rows.Add($refRow)\n$newRow2.Cells.Item(1).Range.Text = \"0.00009\"\n$newRow1 = $tbl.Rows.Add($refRow)\n$newRow1.Cells.Item(1).Range.Text = \"0.00008\"\n\n# --- collapse the three multi-run/tab-separated summary rows near\n#     the end of the table into single plain values\n#     (1-based items 44, 45, 46) ---\n$tbl.Cell(44, 1).Range.Text = \"100\"\n$tbl.Cell(45, 1).Range.Text = \"0\"\n$tbl.Cell(46, 1).Range.Text = \"47\"\n"}
